$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-10 Wednesday" "2024-04-11 Thursday"

Replace-Text "483×9=" "538×3="
Replace-Text "739×8=" "722×7="
Replace-Text "742×2=" "828×5="
Replace-Text "932×5=" "790×2="
Replace-Text "894×3=" "193×3="
Replace-Text "878×9=" "165×7="
Replace-Text "502×2=" "987×8="
Replace-Text "318×2=" "237×9="
Replace-Text "352×9=" "247×6="
Replace-Text "420×8=" "840×7="
Replace-Text "459×6=" "238×4="
Replace-Text "762×4=" "537×3="
Replace-Text "858×8=" "271×3="
Replace-Text "950×3=" "919×6="
Replace-Text "263×8=" "651×4="
Replace-Text "109×5=" "507×2="
Replace-Text "312×2=" "681×8="
Replace-Text "490×5=" "857×2="
Replace-Text "329×3=" "877×5="
Replace-Text "724×2=" "801×6="
Replace-Text "155×2=" "541×9="
Replace-Text "303×6=" "659×2="
Replace-Text "125×3=" "275×6="
Replace-Text "407×6=" "558×5="
Replace-Text "411×9=" "435×7="
